$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$tblShape = $s.Shapes.Item(1)
$tbl = $tblShape.Table

# Remove the row that reads "Lightweight images" / "Heavyweight images"
for ($r = $tbl.Rows.Count; $r -ge 1; $r--) {
    $row = $tbl.Rows.Item($r)
    if ($row.Cells.Item(1).Shape.TextFrame.TextRange.Text -eq "Lightweight images") {
        $row.Delete()
        break
    }
}
